# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 11:22"

# Row 8 - Alemania
$ws.Range("B8").Value = 150729
$ws.Range("C8").Value = 81
$ws.Range("E8").Value = 42114

# Row 15 - Belgica
$ws.Range("B15").Value = 42797
$ws.Range("C15").Value = 908
$ws.Range("D15").Value = 9800
$ws.Range("E15").Value = 26507
$ws.Range("F15").Value = 993
$ws.Range("G15").Value = 228
$ws.Range("H15").Value = 6490

# Row 38 - Indonesia
$ws.Range("B38").Value = 7775
$ws.Range("C38").Value = 357
$ws.Range("D38").Value = 960
$ws.Range("E38").Value = 6168
$ws.Range("G38").Value = 12
$ws.Range("H38").Value = 647

# Row 85 - Hong Kong
$ws.Range("B85").Value = 1036
$ws.Range("C85").Value = 2
$ws.Range("D85").Value = 699
$ws.Range("E85").Value = 333
$ws.Range("F85").Value = 9

# Row 165 - Nepal
$ws.Range("B165").Value = 45
$ws.Range("C165").Value = 0
$ws.Range("E165").Value = 38

# Row 187 - Namibia
$ws.Range("D187").Value = 7
$ws.Range("E187").Value = 9
